$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.448.30'
$ws.Range('E2').Value = '  -5.28%  '
$ws.Range('D3').Value = '3.206.66'
$ws.Range('E3').Value = '  -8.67%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '585.37'
$ws.Range('E5').Value = '  -3.71%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '149.55'
$ws.Range('E6').Value = '  -13.37%  '
$ws.Range('E7').Value = '  +0.21%  '
$ws.Range('D8').Value = '3.199.13'
$ws.Range('E8').Value = '  -8.77%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.543'
$ws.Range('E9').Value = '  -10.62%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.170'
$ws.Range('E10').Value = '  -12.78%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '6.37'
$ws.Range('E11').Value = '  -10.63%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.498'
$ws.Range('E12').Value = '  -15.20%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '38.54'
$ws.Range('E13').Value = '  -16.79%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.0000242'
$ws.Range('E14').Value = '  -12.32%  '
$ws.Range('D15').Value = '3.735.63'
$ws.Range('E15').Value = '  -8.39%  '
$ws.Range('D16').Value = '66.575.83'
$ws.Range('E16').Value = '  -5.15%  '
$ws.Range('D17').Value = '3.220.04'
$ws.Range('E17').Value = '  -8.60%  '
$ws.Range('E18').Value = '  -6.19%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '534.16'
$ws.Range('E19').Value = '  -12.90%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '7.20'
$ws.Range('E20').Value = '  -14.53%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '15.08'
$ws.Range('E21').Value = '  -15.00%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.759'
$ws.Range('E22').Value = '  -13.85%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '7.78'
$ws.Range('E23').Value = '  -13.51%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '13.62'
$ws.Range('E24').Value = '  -12.47%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '84.91'
$ws.Range('E25').Value = '  -14.36%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  +0.16%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '3.15'
$ws.Range('E27').Value = '  -15.89%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.17'
$ws.Range('E28').Value = '  -15.67%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '7.99'
$ws.Range('E29').Value = '  -12.24%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '29.13'
$ws.Range('E30').Value = '  -13.48%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '2.56'
$ws.Range('E31').Value = '  -14.12%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.13'
$ws.Range('E32').Value = '  -14.27%  '
$ws.Range('B33').Value = 'Bittensor'
$ws.Range('C33').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '530.71'
$ws.Range('E33').Value = '  -15.08%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '6.49'
$ws.Range('E34').Value = '  -19.82%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '5.70'
$ws.Range('E35').Value = '  -16.30%  '
$ws.Range('E36').Value = '  +0.04%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '53.18'
$ws.Range('E37').Value = '  -6.43%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.0859'
$ws.Range('E38').Value = '  -14.28%  '
$ws.Range('B39').Value = 'Cosmos'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '9.23'
$ws.Range('E39').Value = '  -14.36%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.0411'
$ws.Range('E40').Value = '  -16.63%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.124'
$ws.Range('E41').Value = '  -14.38%  '
$ws.Range('D42').Value = '2.903.97'
$ws.Range('E42').Value = '  -13.60%  '
$ws.Range('E43').Value = '  -24.85%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.261'
$ws.Range('E44').Value = '  -15.81%  '
$ws.Range('D45').Value = '0.0₃0582'
$ws.Range('E45').Value = '  -21.48%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.37'
$ws.Range('E46').Value = '  -18.22%  '
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '25.92'
$ws.Range('E48').Value = '  -19.21%  '
$ws.Range('B49').Value = 'Fetch.AI'
$ws.Range('C49').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.10'
$ws.Range('E49').Value = '  -17.73%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.114'
$ws.Range('E50').Value = '  -12.80%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '122.14'
$ws.Range('E51').Value = '  -8.31%  '
